$d = $word.ActiveDocument

# 1. Update the cached "Last updated, Date" field result text.
$d.Content.Find.Execute("April 30, 2015", $true, $false, $false, $false, $false, $true, 1, $false, "April 25, 2020", 2) | Out-Null

# 2. Remove the stray "_GoBack" bookmark that previously sat right after the
#    word "Manual" (left over from an earlier edit session).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3. Bump the Revit SDK folder version referenced in the sample App.config
#    snippet from 2013 to 2019 (only the lone digit run changes).
$sdkRange = $d.Content
$sdkRange.Find.Execute("Revit SDK 2013", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$digitRange = $d.Range($sdkRange.End - 1, $sdkRange.End)
$digitRange.Text = "9"

# 4. Word re-drops a "_GoBack" bookmark at the location of the most recent
#    edit - recreate it right after the digit we just changed.
$goBackRange = $d.Range($digitRange.End, $digitRange.End)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
